$wb = $excel.ActiveWorkbook

# Rename "Sheet1" to "Data"
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Name = "Data"

# "Data" (formerly Sheet1) becomes the active/selected tab instead of "Legend"
$sheet1.Activate()
